$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.976.36"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "1.795.18"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'307.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.4186"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("D8").Value = "'0.3556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.06%  "
$ws.Range("D9").Value = "'0.07083"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.75%  "
$ws.Range("D10").Value = "'0.8440"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("D11").Value = "'20.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").Value = "1.794.10"
$ws.Range("E12").Value = "  -4.70%  "
$ws.Range("D13").Value = "'5.279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").Value = "'6.338"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.97%  "
$ws.Range("D15").Value = "'0.06750"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "'79.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "'0.000008654"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.44%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'15.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.79%  "
$ws.Range("D21").Value = "27.059.74"
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").Value = "'5.047"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").Value = "'10.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "2.014.76"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").Value = "'1.936"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "'152.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").Value = "'18.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "'4.991"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.24%  "
$ws.Range("D29").Value = "'112.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").Value = "'1.641"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.67%  "
$ws.Range("D31").Value = "'0.08919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'0.7152"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.29%  "
$ws.Range("D33").Value = "'2.856"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("D34").Value = "'4.289"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.06%  "
$ws.Range("D35").Value = "'1.003"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "'1.069"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.48%  "
$ws.Range("D37").Value = "'1.073"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("D38").Value = "'0.01899"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").Value = "'0.05093"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.37%  "
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").Value = "'0.4935"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").Value = "'2.570"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.49%  "
$ws.Range("D43").Value = "'5.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.41%  "
$ws.Range("D44").Value = "'8.013"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.51%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'104.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").Value = "'1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").Value = "'0.06297"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.14%  "
$ws.Range("D49").Value = "'0.4505"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.04%  "
$ws.Range("D50").Value = "'1.589"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("D51").Value = "'61.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.97%  "
